# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 40 (pushing the existing rows 40-81
# down to 41-82) and populate the new row with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 40..81 down to 41..82.
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new weekly record.
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = 'Vega Monumental Concepción'
$ws.Range("C40").Value = 'Bíobío'
$ws.Range("D40").Value = 44539
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112032
$ws.Range("G40").Value = 'Zapallo italiano'
$ws.Range("H40").Value = 'Sin especificar'
$ws.Range("I40").Value = 'Primera'
$ws.Range("J40").Value = 270
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = 5222
$ws.Range("N40").Value = '$/caja 60 unidades'
$ws.Range("O40").Value = "Región de O'Higgins"
$ws.Range("P40").Value = 87
$ws.Range("Q40").Value = 60
$ws.Range("R40").Value = 'Hortaliza'
